$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 24
$ws.Range("I11").Value = 24
$ws.Range("K11").Value = 24
$ws.Range("M11").Value = 116
$ws.Range("H70").Value = 3043.1428
$ws.Range("I70").Value = 3350.5
$ws.Range("J70").Value = 2633.3333
$ws.Range("K70").Value = 10051.5
$ws.Range("L70").Value = 7899.999899999999
$ws.Range("M70").Value = -9781.5
$ws.Range("N70").Value = -8439.999899999999
$ws.Range("H73").Value = 3043.1428
$ws.Range("I73").Value = 3350.5
$ws.Range("J73").Value = 2633.3333
$ws.Range("K73").Value = 10051.5
$ws.Range("L73").Value = 7899.999899999999
$ws.Range("M73").Value = -9115.5
$ws.Range("N73").Value = -9771.999899999999
$ws.Range("H127").Value = 1921.0526
$ws.Range("J127").Value = 1921.0526
$ws.Range("L127").Value = 5763.1578
$ws.Range("N127").Value = -15683.1578
$ws.Range("H137").Value = 1713.8889
$ws.Range("I137").Value = 1358.8636
$ws.Range("J137").Value = 3276
$ws.Range("K137").Value = 4076.5908
$ws.Range("L137").Value = 9828
$ws.Range("M137").Value = -1526.5908
$ws.Range("N137").Value = -14928
$ws.Range("H138").Value = 3879.8518
$ws.Range("I138").Value = 805.8333
$ws.Range("J138").Value = 7722.375
$ws.Range("K138").Value = 2417.4999
$ws.Range("L138").Value = 23167.125
$ws.Range("M138").Value = 2722.5001
$ws.Range("N138").Value = -33447.125

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 8875
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 7750
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 7750
$ws.Range("M3").Value = -9885
$ws.Range("N3").Value = -7980
$ws.Range("H22").Value = 3733.3333
$ws.Range("I22").Value = 1990
$ws.Range("K22").Value = 1990
$ws.Range("M22").Value = -1691
$ws.Range("H32").Value = 4676.8687
$ws.Range("I32").Value = 3420.7659
$ws.Range("J32").Value = 8893.786
$ws.Range("K32").Value = 3420.7659
$ws.Range("L32").Value = 8893.786
$ws.Range("M32").Value = -3133.7659
$ws.Range("N32").Value = -9467.786

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 83335090
$ws.Range("I86").Value = 83335090
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 83335090
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -83333967
$ws.Range("H89").Value = 83335090
$ws.Range("I89").Value = 83335090
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 416675450
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -416669834
$ws.Range("H94").Value = 1617.5333
$ws.Range("I94").Value = 1420.75
$ws.Range("J94").Value = 2101.923
$ws.Range("K94").Value = 1420.75
$ws.Range("L94").Value = 2101.923
$ws.Range("M94").Value = -969.75
$ws.Range("N94").Value = -3003.923

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9038.706
$ws.Range("I31").Value = 2302.875
$ws.Range("K31").Value = 2302.875
$ws.Range("M31").Value = -2007.875
$ws.Range("H34").Value = 9038.706
$ws.Range("I34").Value = 2302.875
$ws.Range("K34").Value = 2302.875
$ws.Range("M34").Value = -2100.875
$ws.Range("H62").Value = 12500
$ws.Range("I62").Value = 20000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 20000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -19376
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 12500
$ws.Range("I65").Value = 20000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 100000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -96880
$ws.Range("N65").Value = -31240

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 287
$ws.Range("I18").Value = 293.33334
$ws.Range("J18").Value = 230
$ws.Range("K18").Value = 880.0000200000001
$ws.Range("L18").Value = 690
$ws.Range("M18").Value = -711.0000200000001
$ws.Range("N18").Value = -1028
$ws.Range("H124").Value = 6394.75
$ws.Range("I124").Value = 1330
$ws.Range("K124").Value = 3990
$ws.Range("M124").Value = 920
$ws.Range("I125").Value = 1450
$ws.Range("J125").Value = 3844.1875
$ws.Range("K125").Value = 4350
$ws.Range("L125").Value = 11532.5625
$ws.Range("M125").Value = 570
$ws.Range("N125").Value = -21372.5625
$ws.Range("H126").Value = 5108
$ws.Range("J126").Value = 6440
$ws.Range("L126").Value = 19320
$ws.Range("N126").Value = -29200
$ws.Range("H129").Value = 1651.421
$ws.Range("I129").Value = 1483.625
$ws.Range("J129").Value = 1773.4546
$ws.Range("K129").Value = 4450.875
$ws.Range("L129").Value = 5320.3638
$ws.Range("M129").Value = 549.125
$ws.Range("N129").Value = -15320.3638
$ws.Range("H130").Value = 7766.6665
$ws.Range("J130").Value = 7766.6665
$ws.Range("L130").Value = 23299.9995
$ws.Range("N130").Value = -33339.99950000001
$ws.Range("H131").Value = 22000854
$ws.Range("I131").Value = 8333787.5
$ws.Range("J131").Value = 26316770
$ws.Range("K131").Value = 25001362.5
$ws.Range("L131").Value = 78950310
$ws.Range("M131").Value = -24996322.5
$ws.Range("N131").Value = -78960390

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5000.4287
$ws.Range("J10").Value = 5000.4287
$ws.Range("L10").Value = 5000.4287
$ws.Range("N10").Value = -5338.4287
$ws.Range("H113").Value = 125001016
$ws.Range("I113").Value = 200000690
$ws.Range("J113").Value = 1566.6666
$ws.Range("K113").Value = 200000690
$ws.Range("L113").Value = 1566.6666
$ws.Range("M113").Value = -199998520
$ws.Range("N113").Value = -5906.6666

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 50002780
$ws.Range("I40").Value = 62502548
$ws.Range("K40").Value = 62502548
$ws.Range("M40").Value = -62502412
$ws.Range("H62").Value = 35000
$ws.Range("I62").Value = 30000
$ws.Range("K62").Value = 30000
$ws.Range("M62").Value = -29376
$ws.Range("H65").Value = 35000
$ws.Range("I65").Value = 30000
$ws.Range("K65").Value = 90000
$ws.Range("M65").Value = -86880

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 17666.666
$ws.Range("J40").Value = 17666.666
$ws.Range("L40").Value = 17666.666
$ws.Range("N40").Value = -17964.666
$ws.Range("H64").Value = 25109.334
$ws.Range("J64").Value = 25109.334
$ws.Range("L64").Value = 25109.334
$ws.Range("N64").Value = -25605.334
$ws.Range("H67").Value = 25109.334
$ws.Range("J67").Value = 25109.334
$ws.Range("L67").Value = 25109.334
$ws.Range("N67").Value = -26825.334
$ws.Range("H69").Value = 21421.375
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 25228.5
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 25228.5
$ws.Range("M69").Value = -9251
$ws.Range("N69").Value = -26726.5
$ws.Range("H72").Value = 21421.375
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 25228.5
$ws.Range("K72").Value = 30000
$ws.Range("L72").Value = 75685.5
$ws.Range("M72").Value = -26256
$ws.Range("N72").Value = -83173.5
$ws.Range("H136").Value = 4832.533
$ws.Range("I136").Value = 5958.8
$ws.Range("J136").Value = 2580
$ws.Range("K136").Value = 17876.4
$ws.Range("L136").Value = 7740
$ws.Range("M136").Value = -15326.4
$ws.Range("N136").Value = -12840
